# Update the ticker-level figures on both report sheets (RUN -> CRWD) and
# refresh the column widths on IncomeReport that Excel re-flowed afterwards.

$wb = $excel.ActiveWorkbook

$wsIncome   = $wb.Worksheets.Item("IncomeReport")
$wsEarnings = $wb.Worksheets.Item("EarningsReport")

# --- IncomeReport (row 2) ---------------------------------------------
$wsIncome.Range("A2").Value = "CRWD"
$wsIncome.Range("B2").Value = "487.83M -> 535.15M -> 580.88M -> 637.37M -> 692.58M"
$wsIncome.Range("C2").Value = 42
$wsIncome.Range("D2").Value = "-0.14 -> -0.21 -> -0.24 -> -0.2 -> 0.0"
$wsIncome.Range("E2").Value = 100
$wsIncome.Range("F2").Value = "159.74M -> 138.25M -> 176.41M -> 212.85M -> 230.93M"
$wsIncome.Range("G2").Value = 45
$wsIncome.Range("H2").Value = "1.61 <- 1.22 <- 2.19 <- 3.79 <- 4.63"
$wsIncome.Range("I2").Value = "12.49 <- 12.06 <- 20.25 <- 25.68 <- 31.10"

# Column widths shifted slightly as a side effect of the content refresh.
$wsIncome.Columns.Item(4).ColumnWidth = 18.83
$wsIncome.Columns.Item(6).ColumnWidth = 25.5
$wsIncome.Columns.Item(9).ColumnWidth = 20.5

# --- EarningsReport (row 2) --------------------------------------------
$wsEarnings.Range("A2").Value = "CRWD"

# B2/C2 hold numeric-looking text ("0.51", "0.57") in the source sheet, so
# force text entry (leading apostrophe) instead of letting Excel coerce
# them to real numbers, then restore the original cell formatting (the
# apostrophe/quote-prefix trick stamps its own number format otherwise).
$wsEarnings.Range("B2").Value = "'0.51"
$wsEarnings.Range("C2").Value = "'0.57"
$wsEarnings.Range("A2").Copy() | Out-Null
$wsEarnings.Range("B2:C2").PasteSpecial(-4122) | Out-Null
$wsEarnings.Range("D2").Value = 11
$wsEarnings.Range("E2").Value = "55900, 402"
$wsEarnings.Range("F2").Value = "35, 35"
$wsEarnings.Range("G2").Value = 52
